# chore: adapt column header formatting to respective input file names (#7)
#
# Renames the "_old" / "_new" header-column suffixes to the concrete
# format-version names "_FV2210" / "_FV2304", turns the header row +
# data range into a proper Excel Table ("Table1"), and freezes the
# header row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the header cells -------------------------------------------------
# Columns A:J used the "_old" suffix, columns L:U used the "_new" suffix
# (column K just holds the literal "diff" header and is left untouched).
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2210"
}

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + "_FV2304"
}

# --- 2) Turn the used range into an Excel Table ("Table1") ---------------------
$tableRange = $ws.Range("A1:U63")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = [Type]::Missing

# --- 3) Freeze the header row ----------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
